$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "67.876.50"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +0.79%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.638.60"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +0.64%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E4").Value = "  +0.04%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "597.95"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +0.18%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "153.66"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.61%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +0.02%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.550"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -0.49%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.638.15"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +0.72%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.135"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +10.64%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.159"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -0.67%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "5.22"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +0.68%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.347"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -0.07%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "27.63"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +0.12%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.0000188"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +3.62%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.119.12"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +0.82%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "67.821.03"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +0.72%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.646.49"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +1.12%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "11.43"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +2.80%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "373.19"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +2.85%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "7.50"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +0.24%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.25"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -0.99%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "4.81"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -1.55%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.06"
$cell.ClearFormats()
$ws.Range("E24").Value = "  -1.99%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "72.16"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +1.65%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +0.01%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "9.96"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -1.36%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.760.30"
$cell.ClearFormats()
$ws.Range("E28").Value = "  -0.03%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.0000104"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +1.76%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -2.67%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "576.73"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -1.25%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.40"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +0.56%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "7.88"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +0.80%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.84"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +0.44%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +0.09%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.126"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +0.30%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.51"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -0.30%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "157.91"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +0.36%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "19.20"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +0.28%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.91"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +5.55%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.369"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +0.19%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "5.36"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +1.76%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.0₆0340"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +18.72%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.63"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +2.35%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "17.12"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +4.75%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +0.04%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "40.25"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -2.23%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "156.22"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -0.12%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "3.69"
$cell.ClearFormats()
$ws.Range("E49").Value = "  -0.96%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "21.94"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +6.39%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.70"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -1.29%  "
